$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C6) from 2023-10-05 (45204) to 2023-10-08 (45207)
for ($row = 2; $row -le 6; $row++) {
    $ws.Range("C$row").Value = 45207
}
